$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new "Firma" header columns (G, H, I) mirroring column F ---
$ws.Range("G7").Value = "Firma"
$ws.Range("H7").Value = "Firma"
$ws.Range("I7").Value = "Firma"

# Copy F7's formatting (bold font + bottom border from the row-7 header
# style) onto the new header cells so they look the same as the rest of
# the header row.
$ws.Range("F7").Copy() | Out-Null
$ws.Range("G7:I7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Size the three new columns to match the source widths ---
$ws.Cells.Item(1, 7).ColumnWidth = 19.833333333333332   # column G (~20.7265625 chars)
$ws.Cells.Item(1, 8).ColumnWidth = 19.0                 # column H (~19.81640625 chars)
$ws.Cells.Item(1, 9).ColumnWidth = 18.666666666666668   # column I (~19.453125 chars)

# --- Update the selected / active cell shown when the sheet is reopened ---
$ws.Range("A7").Select() | Out-Null
